$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range('D2').Value = '28.417.85'
$ws.Range('E2').Value = '  +0.24%  '

$ws.Range('D3').Value = '1.818.03'
$ws.Range('E3').Value = '  -0.32%  '

$ws.Range('E4').Value = '  +0.09%  '

Set-TextValue 'D5' '315.67'
$ws.Range('E5').Value = '  -0.50%  '

$ws.Range('E6').Value = '  +0.02%  '

Set-TextValue 'D7' '0.5115'
$ws.Range('E7').Value = '  -4.04%  '

Set-TextValue 'D8' '0.3956'
$ws.Range('E8').Value = '  -2.18%  '

Set-TextValue 'D9' '0.08149'
$ws.Range('E9').Value = '  +7.33%  '

$ws.Range('E10').Value = '  -0.40%  '

$ws.Range('E11').Value = '  +0.16%  '

$ws.Range('E12').Value = '  +0.85%  '

Set-TextValue 'D13' '6.279'
$ws.Range('E13').Value = '  -0.58%  '

$ws.Range('E14').Value = '  +0.06%  '

Set-TextValue 'D15' '7.509'
$ws.Range('E15').Value = '  -1.33%  '

$ws.Range('D16').Value = '1.816.71'
$ws.Range('E16').Value = '  -0.56%  '

Set-TextValue 'D17' '0.00001135'
$ws.Range('E17').Value = '  +5.78%  '

Set-TextValue 'D18' '92.62'
$ws.Range('E18').Value = '  +3.66%  '

Set-TextValue 'D19' '0.06637'
$ws.Range('E19').Value = '  +0.56%  '

Set-TextValue 'D20' '17.69'
$ws.Range('E20').Value = '  +0.22%  '

$ws.Range('E21').Value = '  +0.03%  '

Set-TextValue 'D22' '6.099'
$ws.Range('E22').Value = '  +0.03%  '

$ws.Range('D23').Value = '28.450.83'
$ws.Range('E23').Value = '  +0.30%  '

$ws.Range('E24').Value = '  +1.00%  '

Set-TextValue 'D25' '2.261'
$ws.Range('E25').Value = '  +2.62%  '

Set-TextValue 'D26' '21.15'
$ws.Range('E26').Value = '  +2.67%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D27' '155.55'
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.027.85'
$ws.Range('E28').Value = '  -0.52%  '

Set-TextValue 'D29' '2.408'
$ws.Range('E29').Value = '  -1.95%  '

Set-TextValue 'D30' '125.99'
$ws.Range('E30').Value = '  +1.78%  '

$ws.Range('E31').Value = '  +0.16%  '

Set-TextValue 'D32' '1.107'
$ws.Range('E32').Value = '  -1.25%  '

Set-TextValue 'D33' '5.785'
$ws.Range('E33').Value = '  +2.39%  '

Set-TextValue 'D34' '3.651'
$ws.Range('E34').Value = '  +0.09%  '

Set-TextValue 'D35' '0.07022'
$ws.Range('E35').Value = '  -5.28%  '

Set-TextValue 'D36' '0.2227'
$ws.Range('E36').Value = '  -0.18%  '

Set-TextValue 'D37' '5.223'
$ws.Range('E37').Value = '  +0.54%  '

Set-TextValue 'D38' '0.02331'
$ws.Range('E38').Value = '  -0.57%  '

Set-TextValue 'D39' '8.824'
$ws.Range('E39').Value = '  -0.81%  '

Set-TextValue 'D40' '0.6282'
$ws.Range('E40').Value = '  +0.52%  '

Set-TextValue 'D41' '11.31'
$ws.Range('E41').Value = '  +0.18%  '

Set-TextValue 'D42' '1.176'
$ws.Range('E42').Value = '  -0.57%  '

Set-TextValue 'D43' '1.000'
$ws.Range('E43').Value = '  +0.01%  '

$ws.Range('E44').Value = '  +0.35%  '

Set-TextValue 'D45' '13.45'
$ws.Range('E45').Value = '  +0.06%  '

Set-TextValue 'D46' '3.740'
$ws.Range('E46').Value = '  +1.17%  '

Set-TextValue 'D47' '0.5925'
$ws.Range('E47').Value = '  +1.48%  '

Set-TextValue 'D48' '124.77'
$ws.Range('E48').Value = '  -0.10%  '

Set-TextValue 'D49' '1.976'
$ws.Range('E49').Value = '  -0.62%  '

Set-TextValue 'D50' '1.186'
$ws.Range('E50').Value = '  -1.23%  '

Set-TextValue 'D51' '0.06889'
$ws.Range('E51').Value = '  -0.06%  '
